$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the active sheet view: scroll to A55 and select W76 (matches the sheetView change in the diff) ---
$win = $excel.ActiveWindow
$win.ScrollRow = 55
$win.ScrollColumn = 1
[void]$ws.Range("W76").Select()

# --- Correct the migrant-count ("C") column and re-flow the dependent D/E/F columns ---
# Row 17
$ws.Range("C17").Value = 19406.97712
$ws.Range("D17").Value = 2585.4999999999986
$ws.Range("E17").Value = 281.99999999999972
$ws.Range("F17").Value = 415783.74875999975

# Row 18
$ws.Range("C18").Value = 17100.190320000002
$ws.Range("D18").Value = 6433.2999999999956
$ws.Range("E18").Value = 210.99999999999983
$ws.Range("F18").Value = 267414.72803999996

# Row 19
$ws.Range("C19").Value = 18039.2
$ws.Range("D19").Value = 9015.3999999999905
$ws.Range("E19").Value = 204.99999999999977
$ws.Range("F19").Value = 229243.99999999988

# Row 20
$ws.Range("C20").Value = 19598.053680000001
$ws.Range("D20").Value = 27639.699999999986
$ws.Range("E20").Value = 219.99999999999972
$ws.Range("F20").Value = 794230.03655999957

# Row 21
$ws.Range("C21").Value = 24531.17148999999
$ws.Range("D21").Value = 9993.7999999999938
$ws.Range("E21").Value = 349.99999999999949
$ws.Range("F21").Value = 135007.86686999991

# Row 47
$ws.Range("C47").Value = 20427.341848341399
$ws.Range("D47").Value = 3380.6948459023511
$ws.Range("E47").Value = 399.82310720475499
$ws.Range("F47").Value = 795019.6038440204

# Row 48
$ws.Range("C48").Value = 20427.341848341399
$ws.Range("D48").Value = 8411.9219308232787
$ws.Range("E48").Value = 299.15842418511818
$ws.Range("F48").Value = 375973.016231196

# Row 49
$ws.Range("C49").Value = 20427.341848341399
$ws.Range("D49").Value = 11788.171074743001
$ws.Range("E49").Value = 290.65154956374028
$ws.Range("F49").Value = 322306.69853012729

# Row 50
$ws.Range("C50").Value = 20427.341848341399
$ws.Range("D50").Value = 27639.69999999999
$ws.Range("E50").Value = 311.91873611718472
$ws.Range("F50").Value = 795019.6038440204

# Row 51
$ws.Range("C51").Value = 24531.17148999999
$ws.Range("D51").Value = 18620.078602856651
$ws.Range("E51").Value = 496.23435291370282
$ws.Range("F51").Value = 189814.95631931321

# Row 77
$ws.Range("C77").Value = 20815.4167831292
$ws.Range("D77").Value = 3187.1120249077162
$ws.Range("E77").Value = 345.11407745492488
$ws.Range("F77").Value = 602086.7667823683

# Row 78
$ws.Range("C78").Value = 20815.4167831292
$ws.Range("D78").Value = 7930.2447456348118
$ws.Range("E78").Value = 258.22365369854322
$ws.Range("F78").Value = 387237.04203390318

# Row 79
$ws.Range("C79").Value = 20815.4167831292
$ws.Range("D79").Value = 11113.165635023401
$ws.Range("E79").Value = 250.8808009867362
$ws.Range("F79").Value = 331962.89940598019

# Row 80
$ws.Range("C80").Value = 20815.4167831292
$ws.Range("D80").Value = 27639.69999999999
$ws.Range("E80").Value = 269.23793276625338
$ws.Range("F80").Value = 1094086.860592877

# Row 81
$ws.Range("C81").Value = 29031.0083701363
$ws.Range("D81").Value = 12134.721466634899
$ws.Range("E81").Value = 419.79132008343811
$ws.Range("F81").Value = 192827.72979110759
